# Penalty Reward System (unfinished) - update forecast dates and zero out MyForecast
$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Helper: write a plain text value into a cell without letting Excel's
# automatic type recognition (e.g. date parsing) change its stored type.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# New Week_Start_Date values: every week shifts forward by one week,
# and a new trailing week (2025-04-27) is appended for row 17.
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $dateCell = $wsForecast.Cells.Item($row, 2)   # column B: Week_Start_Date
    Set-TextValue $dateCell $newDates[$i]

    $forecastCell = $wsForecast.Cells.Item($row, 4)  # column D: MyForecast
    $forecastCell.Value = 0
}

# Update the Summary sheet metrics
Set-TextValue $wsSummary.Range("B2")  "2022-12-25 to 2025-01-05"
Set-TextValue $wsSummary.Range("B9")  "4"
Set-TextValue $wsSummary.Range("B10") "3"
Set-TextValue $wsSummary.Range("B11") "2"
Set-TextValue $wsSummary.Range("B12") "0"
Set-TextValue $wsSummary.Range("B13") "2025-01-12"
Set-TextValue $wsSummary.Range("B14") "0"
Set-TextValue $wsSummary.Range("B15") "2025-02-23"
